# Adds a "Ready for handoff" row for file 9665fb5b-... across the
# Overview / zh-cn / de-de sheets, mirroring the existing d287a20c-... row.

$wb = $excel.ActiveWorkbook

$mdName    = '9665fb5b-2ec3-4674-81ee-d97f77b8a522ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$mdPath    = 'e2e\9665fb5b-2ec3-4674-81ee-d97f77b8a522ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$ready     = 'Ready for handoff'
$dt3931    = '2016-08-14 02:39:31'
$zhcnXlf   = '9665fb5b-2ec3-4674-81ee-d97f77b8a522ooooooooooooooooooooooooooooooooooooooooooo.05d977c2378f3d7a472062e86ce0ea3bfb0c182f.zh-cn.xlf'
$dt3923    = '2016-08-14 02:39:23'
$dedeXlf   = '9665fb5b-2ec3-4674-81ee-d97f77b8a522ooooooooooooooooooooooooooooooooooooooooooo.05d977c2378f3d7a472062e86ce0ea3bfb0c182f.de-de.xlf'
$hlUrl     = 'https://github.com/OpenLocalizationTestOrg/oltest/blob/16f246addb3bb5be8be61093823d85677b330a55/e2e/9665fb5b-2ec3-4674-81ee-d97f77b8a522ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'

# ---------------------------------------------------------------
# Sheet "Overview" (sheet index 1): columns A..G, new row 3
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $mdName
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = $ready
$wsOverview.Range("F3").Value = $ready
$wsOverview.Range("G3").Value = $dt3931

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $hlUrl, "", "", $mdPath)

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------
# Sheet "zh-cn" (sheet index 2): columns A..P, new row 3
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A3").Value = $mdName
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = $ready
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = $zhcnXlf
$wsZhCn.Range("H3").Value = $dt3923
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("O3").Value = "False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $hlUrl, "", "", $mdName)

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------
# Sheet "de-de" (sheet index 3): columns A..P, new row 3
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A3").Value = $mdName
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = $ready
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = $dedeXlf
$wsDeDe.Range("H3").Value = $dt3931
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("O3").Value = "False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $hlUrl, "", "", $mdName)

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))

# ---------------------------------------------------------------
# Column width tweaks (status / zh-cn / de-de columns widened to
# fit the new "Ready for handoff" text)
# ---------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797
$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
